# Updates cryptos list data (coin name/link/price/volume columns)
# to match refreshed scrape results, per commit message:
# "Updated cryptos list on Mon Mar 18 08:34:00 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '67.536.39'
$ws.Range('E2').Value = '  +2.65%  '
# Row 3
$ws.Range('D3').Value = '3.570.90'
$ws.Range('E3').Value = '  +1.76%  '
# Row 4
$ws.Range('E4').Value = '  +0.45%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '199.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +8.01%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '567.62'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.53%  '
# Row 7
$ws.Range('B7').Value = 'LidoStakedEther'
$ws.Range('C7').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D7').Value = '3.564.60'
$ws.Range('E7').Value = '  +1.65%  '
# Row 8
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.612'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.49%  '
# Row 9
$ws.Range('B9').Value = 'USDC'
$ws.Range('C9').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.36%  '
# Row 10
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.673'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.88%  '
# Row 11
$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '59.96'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +14.75%  '
# Row 12
$ws.Range('B12').Value = 'Dogecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.147'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.32%  '
# Row 13
$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000280'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +10.48%  '
# Row 14
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.22'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.37%  '
# Row 15
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '4.167.11'
$ws.Range('E15').Value = '  +2.43%  '
# Row 16
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.591.60'
$ws.Range('E16').Value = '  +2.95%  '
# Row 17
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.126'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.32%  '
# Row 18
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.01'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +6.47%  '
# Row 19
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '67.437.87'
$ws.Range('E19').Value = '  +2.68%  '
# Row 20
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.17'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.78%  '
# Row 21
$ws.Range('B21').Value = 'Polygon'
$ws.Range('C21').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.06'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.27%  '
# Row 22
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '401.20'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.05%  '
# Row 23
$ws.Range('B23').Value = 'RenderToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.82'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +18.63%  '
# Row 24
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.16'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.53%  '
# Row 25
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.49'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.46%  '
# Row 26
$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.88'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.15%  '
# Row 27
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.88'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +11.36%  '
# Row 28
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.34'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.32%  '
# Row 29
$ws.Range('B29').Value = 'LEO'
$ws.Range('C29').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.10'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.73%  '
# Row 30
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.16'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.26%  '
# Row 31
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.66'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.17%  '
# Row 32
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '31.31'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.92%  '
# Row 33
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '670.61'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +9.95%  '
# Row 34
$ws.Range('B34').Value = 'Cosmos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '12.06'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.24%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '63.33'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.24%  '
# Row 36
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.113'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.21%  '
# Row 37
$ws.Range('B37').Value = 'InjectiveProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '41.26'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.69%  '
# Row 38
$ws.Range('B38').Value = 'TheGraph'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.407'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.24%  '
# Row 39
$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.13%  '
# Row 40
$ws.Range('D40').Value = '0.0₃0755'
$ws.Range('E40').Value = '  +4.75%  '
# Row 41
$ws.Range('B41').Value = 'ThetaToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.17'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +14.05%  '
# Row 42
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '3.172.72'
$ws.Range('E42').Value = '  +5.84%  '
# Row 43
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.132'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.28%  '
# Row 44
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.32%  '
# Row 45
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.68'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +8.52%  '
# Row 46
$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.80'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +14.21%  '
# Row 47
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.76'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +18.61%  '
# Row 48
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0409'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.28%  '
# Row 49
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.130'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.90%  '
# Row 50
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.09'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.22%  '
# Row 51
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.57'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.99%  '
